$wb = $excel.ActiveWorkbook

$wsQuery = $wb.Worksheets.Item("Query")
$wsPrepared = $wb.Worksheets.Item("Prepared")

# Create the new "LessGreater" sheet as a copy of "Query" (placed after "Prepared"),
# so it inherits the identical column widths / row heights / cell styles.
$wsQuery.Copy($null, $wsPrepared)
$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "LessGreater"

# Update the new sheet's query cell (A2) to the "not equal" flavored query.
$newQuery = "<jt:forEach items=""`${jdbc.execQuery('SELECT * FROM employee WHERE first_name <> \'Randy\'')}"" var=""employee"" >`${employee.first_name}"
$wsNew.Range("A2").Value = $newQuery

# Update the selections left on the original two sheets.
$wsQuery.Range("A1:G2").Select()
$wsPrepared.Range("A3").Select()

# Leave the new sheet as the active tab/selection.
$wsNew.Activate()
$wsNew.Range("A1").Select()
